$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row translations (Swedish -> English)
$ws.Range("A1").Value = "Name"
$ws.Range("B1").Value = "Age"
$ws.Range("C1").Value = "Email"
$ws.Range("D1").Value = "Mobile"
$ws.Range("E1").Value = "Pet Name"
$ws.Range("F1").Value = "Breed"

# Summary labels translations (Swedish -> English)
$ws.Range("A48").Value = "Average age:"
$ws.Range("A49").Value = "Median age:"
$ws.Range("A50").Value = "Youngest:"
$ws.Range("A51").Value = "Oldest:"
